$d = $word.ActiveDocument

# --- "Updated typos after teaching" ------------------------------------
# The example HTTP request path changes from '/?q=node.js' to
# '/search?q=node.js' (the word "search" is inserted right after the
# leading slash). The document's "_GoBack" bookmark — which currently
# sits right after "check out the docs" — ends up relocated to sit
# between the newly-typed "search" and the rest of the query string,
# reflecting where the author's cursor was left after the edit.

# Locate the '/?q=node.js' inline-code span.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "/?q=node.js"
$find.Forward = $true
$find.Wrap = 0
$null = $find.Execute()
$urlRange = $find.Parent.Duplicate
$slashEnd = $urlRange.Start + 1   # position right after the leading '/'

# Drop the old "_GoBack" bookmark from its current position.
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

# Type "search" right after the '/'.
$insertion = $d.Range($slashEnd, $slashEnd)
$insertion.InsertBefore("search")

# Force the newly-typed word into its own run (distinct from the
# surrounding, identically-styled text) the way a live edit session
# would, by nudging and then reverting a direct character attribute.
$searchRange = $d.Range($slashEnd, $slashEnd + 6)
$searchRange.Bold = 1
$searchRange.Bold = 0

# Re-create "_GoBack" at the new cursor position: right after "search",
# just before "?q=node.js'".
$newPos = $slashEnd + 6
$bmRange = $d.Range($newPos, $newPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
